# Update "想去人数" (attendance interest count) figures in column F
# across the sheets, as produced by the gh-pages data refresh at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F3").Value  = 7868
$ws1.Range("F5").Value  = 935
$ws1.Range("F6").Value  = 290
$ws1.Range("F8").Value  = 603
$ws1.Range("F9").Value  = 91
$ws1.Range("F13").Value = 3143
$ws1.Range("F14").Value = 201
$ws1.Range("F15").Value = 95
$ws1.Range("F16").Value = 737
$ws1.Range("F17").Value = 753
$ws1.Range("F18").Value = 48
$ws1.Range("F19").Value = 456
$ws1.Range("F21").Value = 246
$ws1.Range("F22").Value = 223
$ws1.Range("F23").Value = 311
$ws1.Range("F24").Value = 287
$ws1.Range("F26").Value = 106
$ws1.Range("F27").Value = 274
$ws1.Range("F28").Value = 18
$ws1.Range("F31").Value = 501
$ws1.Range("F32").Value = 515
$ws1.Range("F33").Value = 20
$ws1.Range("F37").Value = 94

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 205

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 205
$ws4.Range("F5").Value  = 7868
$ws4.Range("F7").Value  = 935
$ws4.Range("F8").Value  = 290
$ws4.Range("F10").Value = 603
$ws4.Range("F11").Value = 91
$ws4.Range("F16").Value = 3143
$ws4.Range("F17").Value = 201
$ws4.Range("F18").Value = 95
$ws4.Range("F20").Value = 737
$ws4.Range("F21").Value = 753
$ws4.Range("F23").Value = 48
$ws4.Range("F24").Value = 456
$ws4.Range("F26").Value = 246
$ws4.Range("F27").Value = 223
$ws4.Range("F28").Value = 311
$ws4.Range("F29").Value = 287
$ws4.Range("F31").Value = 106
$ws4.Range("F32").Value = 274
$ws4.Range("F33").Value = 18
$ws4.Range("F36").Value = 501
$ws4.Range("F37").Value = 515
$ws4.Range("F38").Value = 20
$ws4.Range("F42").Value = 94
